# Saldo.xlsx update ("Add files via upload"):
#   - Insert a new record (DILSON / 004472404 / 9471.13) right above the
#     HELIO (004363250) row.
#   - Remove the MIRELLA (003553997 / 1100) record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# HELIO / 004363250 currently sits on row 13 (row 1 is the header).
# Push it (and everything below) down one row so we can drop the new
# DILSON record into the freed-up row 13.
$ws.Rows.Item(13).Insert()

# Account numbers are text (leading zeros must survive), so force the
# cell to text before typing the value, then drop the formatting so the
# cell ends up styled just like its neighbours.
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "004472404"
$ws.Range("A13").ClearFormats()

$ws.Range("B13").Value = "DILSON"
$ws.Range("C13").Value = 9471.13

# After the insert above, the MIRELLA / 003553997 row (originally row 16)
# is now row 17 - remove it entirely.
$ws.Rows.Item(17).Delete()
